# "Logged Week 15 and simulated Week 16"
#
# A new player, G.Gilbert, is added to the roster. He is inserted into the
# lineup right after K.Allen / before A.Gibson - i.e. a new column is
# inserted before column E on both the "Rushing" and "Receiving" sheets,
# pushing the existing players (and the trailing placeholder data) one
# column to the right. The new column gets the player's header name in
# row 1 and the same "n" placeholder value used by every other player in
# row 2.
#
# Apply the same change to every worksheet in the workbook (both the
# "Rushing" and "Receiving" tabs use an identical layout).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new column before column E, shifting E:T -> F:U.
    $ws.Columns("E:E").Insert()

    # New player header (row 1) and placeholder stat value (row 2).
    $ws.Range("E1").Value = "G.Gilbert"
    $ws.Range("E2").Value = "n"
}
